$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.794476509094238
$ws.Range("B1").Value = 1.782221436500549
$ws.Range("C1").Value = 7.684327125549316
$ws.Range("D1").Value = 1.0156090259552
$ws.Range("E1").Value = 0.4105211794376373
